$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.079.88"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").Value = "2.105.56"
$ws.Range("E3").Value = "  +10.33%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'253.23"
$ws.Range("E5").Value = "  +1.63%  "

$ws.Range("D6").Value = "'0.662"
$ws.Range("E6").Value = "  -4.71%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  +5.54%  "

$ws.Range("D9").Value = "'60.78"
$ws.Range("E9").Value = "  +5.15%  "

$ws.Range("E10").Value = "  +1.25%  "

$ws.Range("D11").Value = "'0.0749"
$ws.Range("E11").Value = "  -1.07%  "

$ws.Range("E12").Value = "  +6.81%  "

$ws.Range("D13").Value = "'14.86"
$ws.Range("E13").Value = "  +0.92%  "

$ws.Range("D14").Value = "2.407.41"
$ws.Range("E14").Value = "  +10.19%  "

$ws.Range("D15").Value = "'0.840"
$ws.Range("E15").Value = "  +3.61%  "

$ws.Range("D16").Value = "2.099.39"
$ws.Range("E16").Value = "  +10.04%  "

$ws.Range("D17").Value = "'5.18"
$ws.Range("E17").Value = "  +1.85%  "

$ws.Range("D18").Value = "37.016.27"
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("D19").Value = "'73.38"
$ws.Range("E19").Value = "  -1.17%  "

$ws.Range("D20").Value = "0.0₃0826"
$ws.Range("E20").Value = "  -3.19%  "

$ws.Range("D21").Value = "'13.34"
$ws.Range("E21").Value = "  -2.25%  "

$ws.Range("D22").Value = "'241.91"
$ws.Range("E22").Value = "  -3.32%  "

$ws.Range("D23").Value = "'5.33"
$ws.Range("E23").Value = "  +3.69%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").Value = "'2.52"
$ws.Range("E25").Value = "  -0.57%  "

$ws.Range("D26").Value = "'170.27"
$ws.Range("E26").Value = "  +1.97%  "

$ws.Range("D27").Value = "'9.45"
$ws.Range("E27").Value = "  +7.80%  "

$ws.Range("D28").Value = "'21.15"
$ws.Range("E28").Value = "  +13.33%  "

$ws.Range("D29").Value = "'2.03"
$ws.Range("E29").Value = "  -7.03%  "

$ws.Range("D30").Value = "'26.63"
$ws.Range("E30").Value = "  +39.73%  "

$ws.Range("E31").Value = "  -4.15%  "

$ws.Range("D32").Value = "'1.07"
$ws.Range("E32").Value = "  +23.49%  "

$ws.Range("E33").Value = "  -1.52%  "

$ws.Range("D34").Value = "'0.0618"
$ws.Range("E34").Value = "  +1.06%  "

$ws.Range("E35").Value = "  +5.83%  "

$ws.Range("E36").Value = "  +21.69%  "

$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("E38").Value = "  -4.23%  "

$ws.Range("D39").Value = "'1.81"
$ws.Range("E39").Value = "  -6.67%  "

$ws.Range("E40").Value = "  -9.24%  "

$ws.Range("E41").Value = "  -0.71%  "

$ws.Range("D42").Value = "'1.17"
$ws.Range("E42").Value = "  +7.53%  "

$ws.Range("E43").Value = "  -5.58%  "

$ws.Range("E44").Value = "  -4.92%  "

$ws.Range("E45").Value = "  -2.26%  "

$ws.Range("D46").Value = "1.351.42"
$ws.Range("E46").Value = "  +0.25%  "

$ws.Range("D47").Value = "'0.0856"
$ws.Range("E47").Value = "  +4.87%  "

$ws.Range("D48").Value = "'7.11"
$ws.Range("E48").Value = "  +10.88%  "

$ws.Range("E49").Value = "  +3.92%  "

$ws.Range("D50").Value = "2.293.26"
$ws.Range("E50").Value = "  +9.81%  "

$ws.Range("D51").Value = "'2.27"
$ws.Range("E51").Value = "  -4.08%  "
